$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; existing rows 16-24 shift down to 17-25.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 45049
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100114007
$ws.Cells.Item(16, 7).Value = "Jengibre"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 300
$ws.Cells.Item(16, 11).Value = 13000
$ws.Cells.Item(16, 12).Value = 14000
$ws.Cells.Item(16, 13).Value = 13500
$ws.Cells.Item(16, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 1038
$ws.Cells.Item(16, 17).Value = 13
$ws.Cells.Item(16, 18).Value = "Hortaliza"
